# Update row 8 (year 2025) metrics in the "metricas_recorrencia_anual" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1076
$ws.Range("D8").Value = 176
$ws.Range("E8").Value = 900
$ws.Range("F8").Value = 7.219031993437245
$ws.Range("G8").Value = 83.64312267657994
$ws.Range("H8").Value = 16.35687732342008
